$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Change 1 (slide 11 - "Register" message body_data format):
#   {id};{password}  ->  {id};{password};{name}
#
# The trailing "}" run of the "Register" paragraph is split into three runs:
#   "};{"  "name"  "}"
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$shape11 = $slide11.Shapes.Item(2)
$tr11 = $shape11.TextFrame.TextRange

$full = $tr11.Text
$afterRegister = $full.IndexOf("Register") + ("Register").Length
# the final "}" of "{id};{password}" sits 16 characters after "Register"
$closeBrace = $tr11.Characters($afterRegister + 16, 1)
$closeBrace.Text = "};{name}"

$full = $tr11.Text
$splitPos = $full.IndexOf("};{name}")
$firstPart = $tr11.Characters($splitPos + 1, 3)
$firstPart.Text = "};{"

$full = $tr11.Text
$namePos = $full.IndexOf("name}")
$nameRun = $tr11.Characters($namePos + 1, 4)
$nameRun.Text = "name"

# ---------------------------------------------------------------------------
# Change 2 (slide 9 - "CreateGroup" message body_data format):
#   {id_group};{nama}; {id_1};{id_2}:..:{id_n}
#   Two adjacent runs ("};{nama" and "}; {id_1};{id_2}:..:{") are merged
#   into a single run "};{nama}; {id_1};{id_2}:..:{"
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$shape9 = $slide9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange

$full9 = $tr9.Text
$mergeStart = $full9.IndexOf("};{nama}; {id_1};{id_2}:..:{")
$mergeLen = ("};{nama}; {id_1};{id_2}:..:{").Length
$mergeRange = $tr9.Characters($mergeStart + 1, $mergeLen)
$mergeRange.Text = "};{nama}; {id_1};{id_2}:..:{"
